$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "RÉCOLEMENT"
$ws.Range("C7").Value = "PLAN DE RÉCOLEMENT"

$ws.Range("B9").Value = "REPÉRAGE"
$ws.Range("C9").Value = "PLAN DE REPÉRAGE"

$ws.Range("B10").Value = "ÉTAT DES LIEUX"

$ws.Range("B16").Value = "COPROPRIÉTÉ"
$ws.Range("C16").Value = "PLAN DE COPROPRIÉTÉ"

$ws.Range("B17").Value = "INTÉRIEUR"

$ws.Range("B32").Value = "ÉPURE"
$ws.Range("C32").Value = "ÉPURE D’IMPLANTATION"

$ws.Range("B19").Value = "FAÇADE"
$ws.Range("C19").Value = "PLAN DE FAÇADE"

$ws.Range("C19").Select()
